$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete year rows (2000, 2002, 2005, 2007) - rows 2 through 5.
# This shifts the existing 2010/2012/2015/2017 rows up to rows 2-5.
$ws.Rows("2:5").Delete()

# Append the new 2020 data as row 6, reusing the same label formatting
# (bold, centered, bordered) as the other year cells in column A.
$ws.Cells.Item(5, 1).Copy()
$ws.Cells.Item(6, 1).PasteSpecial(-4122)
$ws.Cells.Item(6, 1).Value = "2020年"

$ws.Cells.Item(6, 2).Value = 3366560.16355687
$ws.Cells.Item(6, 3).Value = 685773033.710874
$ws.Cells.Item(6, 4).Value = 221055004.423882
$ws.Cells.Item(6, 6).Value = 3023440277.96893
$ws.Cells.Item(6, 10).Value = 221880841.059148
$ws.Cells.Item(6, 11).Value = 100482984.284874
$ws.Cells.Item(6, 12).Value = 99984358.6309073
$ws.Cells.Item(6, 13).Value = 108193778.357945
$ws.Cells.Item(6, 15).Value = 468986.096341701
$ws.Cells.Item(6, 16).Value = 4215024.12649486
$ws.Cells.Item(6, 18).Value = 3925831.46539935
$ws.Cells.Item(6, 19).Value = 452326731.017457
